$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("configurations")

# Helper: force a cell's value to be stored as literal text (shared string),
# never auto-coerced into a boolean/number by Excel's input parser.
# We do this by writing a text-returning formula into a scratch cell,
# copying it, and pasting-special (values only) into the destination -
# Excel keeps pasted formula results as plain text, not re-parsed input.
function Set-LiteralText {
    param($sheet, $cellRef, [string]$text)
    $scratch = $sheet.Range("Z1")
    $scratch.Formula = '="' + $text + '"'
    $scratch.Copy()
    $sheet.Range($cellRef).PasteSpecial(-4163) | Out-Null  # xlPasteValues
    $scratch.Clear()
}

# 1) "Optimize_DMO_name" row (B3): DMO_optimized_test -> False (placeholder text, still the
#    same shared-string slot as "DMO_optimized_test" since B3 is its only user) then reused
#    by B10/B12, and finally overwritten to its real new value "Optimized_DMO".
Set-LiteralText $ws "B3" "False"

# 2) Copy that literal "False" text into the two rows that should now read False
#    (report_dependencies / B10, report_add_optimize / B12) instead of True.
$ws.Range("B3").Copy()
$ws.Range("B10").PasteSpecial(-4163) | Out-Null
$ws.Range("B12").PasteSpecial(-4163) | Out-Null

# 3) Now give B3 its real new value.
Set-LiteralText $ws "B3" "Optimized_DMO"

$ws.Range("D1").Clear()

# Update the sheet view: new zoom level and selected cell.
$ws.Activate()
$excel.ActiveWindow.Zoom = 189
$ws.Range("B13").Select()
